$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1250312.5
$ws.Range("I9").Value = 2500525.2
$ws.Range("J9").Value = 99.75
$ws.Range("K9").Value = 2500525.2
$ws.Range("L9").Value = 99.75
$ws.Range("M9").Value = -2500356.2
$ws.Range("N9").Value = -437.75
$ws.Range("H53").Value = 218.76923
$ws.Range("I53").Value = 70.25
$ws.Range("J53").Value = 456.4
$ws.Range("K53").Value = 70.25
$ws.Range("L53").Value = 456.4
$ws.Range("M53").Value = 566.75
$ws.Range("N53").Value = -1730.4
$ws.Range("H62").Value = 97741.586
$ws.Range("I62").Value = 126011
$ws.Range("J62").Value = 12933.333
$ws.Range("K62").Value = 126011
$ws.Range("L62").Value = 12933.333
$ws.Range("M62").Value = -125387
$ws.Range("N62").Value = -14181.333
$ws.Range("H65").Value = 97741.586
$ws.Range("I65").Value = 126011
$ws.Range("J65").Value = 12933.333
$ws.Range("K65").Value = 630055
$ws.Range("L65").Value = 64666.665
$ws.Range("M65").Value = -626935
$ws.Range("N65").Value = -70906.66500000001
$ws.Range("H132").Value = 402446.28
$ws.Range("I132").Value = 2762.647
$ws.Range("J132").Value = 1251774
$ws.Range("K132").Value = 8287.940999999999
$ws.Range("L132").Value = 3755322
$ws.Range("M132").Value = -5757.940999999999
$ws.Range("N132").Value = -3760382
$ws.Range("H137").Value = 1066.3334
$ws.Range("I137").Value = 943.9167
$ws.Range("K137").Value = 2831.7501
$ws.Range("M137").Value = -281.7501000000002
$ws.Range("H138").Value = 2805.84
$ws.Range("I138").Value = 808.7778
$ws.Range("J138").Value = 3244.2195
$ws.Range("K138").Value = 2426.3334
$ws.Range("L138").Value = 9732.658500000001
$ws.Range("M138").Value = 2713.6666
$ws.Range("N138").Value = -20012.6585

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10340
$ws.Range("H20").Value = 10000
$ws.Range("J20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("N20").Value = -10540
$ws.Range("H45").Value = 1926
$ws.Range("I45").Value = 1687.3334
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 1687.3334
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -1310.3334
$ws.Range("N45").Value = -3754
$ws.Range("H61").Value = 1411.5294
$ws.Range("I61").Value = 1265.4073
$ws.Range("J61").Value = 1975.1428
$ws.Range("K61").Value = 1265.4073
$ws.Range("L61").Value = 1975.1428
$ws.Range("M61").Value = -1053.4073
$ws.Range("N61").Value = -2399.1428
$ws.Range("H74").Value = 913.4
$ws.Range("I74").Value = 869.8158
$ws.Range("J74").Value = 1150
$ws.Range("K74").Value = 869.8158
$ws.Range("L74").Value = 1150
$ws.Range("M74").Value = 4.184200000000033
$ws.Range("N74").Value = -2898
$ws.Range("H77").Value = 913.4
$ws.Range("I77").Value = 869.8158
$ws.Range("J77").Value = 1150
$ws.Range("K77").Value = 4349.079
$ws.Range("L77").Value = 5750
$ws.Range("M77").Value = 18.92100000000028
$ws.Range("N77").Value = -14486
$ws.Range("H97").Value = 648.5263
$ws.Range("I97").Value = 628.25
$ws.Range("J97").Value = 756.6667
$ws.Range("K97").Value = 628.25
$ws.Range("L97").Value = 756.6667
$ws.Range("M97").Value = -132.25
$ws.Range("N97").Value = -1748.6667
$ws.Range("H110").Value = 558.0833
$ws.Range("I110").Value = 558.0833
$ws.Range("K110").Value = 558.0833
$ws.Range("M110").Value = 1486.9167
$ws.Range("H122").Value = 1899
$ws.Range("I122").Value = 1899
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5697
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3247
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1906.409
$ws.Range("I132").Value = 1276.3334
$ws.Range("J132").Value = 2662.5
$ws.Range("K132").Value = 3829.0002
$ws.Range("L132").Value = 7987.5
$ws.Range("M132").Value = -1299.0002
$ws.Range("N132").Value = -13047.5
$ws.Range("H136").Value = 1411.5294
$ws.Range("I136").Value = 1265.4073
$ws.Range("J136").Value = 1975.1428
$ws.Range("K136").Value = 3796.2219
$ws.Range("L136").Value = 5925.428400000001
$ws.Range("M136").Value = -1246.2219
$ws.Range("N136").Value = -11025.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 896
$ws.Range("I94").Value = 602.5714
$ws.Range("K94").Value = 602.5714
$ws.Range("M94").Value = -151.5714
$ws.Range("H107").Value = 8555
$ws.Range("I107").Value = 750.38464
$ws.Range("K107").Value = 750.38464
$ws.Range("M107").Value = 1169.61536
$ws.Range("H134").Value = 27877.39
$ws.Range("I134").Value = 4022.5715
$ws.Range("J134").Value = 52924.95
$ws.Range("K134").Value = 12067.7145
$ws.Range("L134").Value = 158774.85
$ws.Range("M134").Value = -9532.7145
$ws.Range("N134").Value = -163844.85

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6290
$ws.Range("I31").Value = 6290
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6290
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -5995
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 6290
$ws.Range("I34").Value = 6290
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 6290
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -6088
$ws.Range("N34").ClearContents()
$ws.Range("H58").Value = 6002.857
$ws.Range("I58").Value = 963.75
$ws.Range("J58").Value = 9103.846
$ws.Range("K58").Value = 963.75
$ws.Range("L58").Value = 9103.846
$ws.Range("M58").Value = -760.75
$ws.Range("N58").Value = -9509.846
$ws.Range("H76").Value = 83337910
$ws.Range("I76").Value = 83337910
$ws.Range("K76").Value = 83337910
$ws.Range("M76").Value = -83337595
$ws.Range("H79").Value = 83337910
$ws.Range("I79").Value = 83337910
$ws.Range("K79").Value = 83337910
$ws.Range("M79").Value = -83336818
$ws.Range("H132").Value = 2859.926
$ws.Range("I132").Value = 2206.8333
$ws.Range("J132").Value = 4166.1113
$ws.Range("K132").Value = 6620.499899999999
$ws.Range("L132").Value = 12498.3339
$ws.Range("M132").Value = -4090.499899999999
$ws.Range("N132").Value = -17558.3339
$ws.Range("H134").Value = 2966.85
$ws.Range("I134").Value = 2209.7856
$ws.Range("J134").Value = 4733.3335
$ws.Range("K134").Value = 6629.3568
$ws.Range("L134").Value = 14200.0005
$ws.Range("M134").Value = -4094.3568
$ws.Range("N134").Value = -19270.0005
$ws.Range("H136").Value = 6002.857
$ws.Range("I136").Value = 963.75
$ws.Range("J136").Value = 9103.846
$ws.Range("K136").Value = 2891.25
$ws.Range("L136").Value = 27311.538
$ws.Range("M136").Value = -341.25
$ws.Range("N136").Value = -32411.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 3655
$ws.Range("J98").Value = 20000
$ws.Range("L98").Value = 60000
$ws.Range("N98").Value = -62996
$ws.Range("H122").Value = 357999.4
$ws.Range("J122").Value = 834600.75
$ws.Range("L122").Value = 7511406.75
$ws.Range("N122").Value = -7516306.75
$ws.Range("H132").Value = 500
$ws.Range("I132").Value = 500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1970
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1700.3334
$ws.Range("I97").Value = 1700.3334
$ws.Range("K97").Value = 1700.3334
$ws.Range("M97").Value = -1204.3334
$ws.Range("H132").Value = 3237.25
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 3983
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 11949
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -17009

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1928
$ws.Range("I7").Value = 1993.1538
$ws.Range("J7").Value = 1833.8889
$ws.Range("K7").Value = 1993.1538
$ws.Range("L7").Value = 1833.8889
$ws.Range("M7").Value = -1881.1538
$ws.Range("N7").Value = -2057.8889
$ws.Range("H22").Value = 1383.591
$ws.Range("I22").Value = 858.3333
$ws.Range("J22").Value = 1580.5625
$ws.Range("K22").Value = 858.3333
$ws.Range("L22").Value = 1580.5625
$ws.Range("M22").Value = -563.3333
$ws.Range("N22").Value = -2170.5625
$ws.Range("H27").Value = 1383.591
$ws.Range("I27").Value = 858.3333
$ws.Range("J27").Value = 1580.5625
$ws.Range("K27").Value = 858.3333
$ws.Range("L27").Value = 1580.5625
$ws.Range("M27").Value = -751.3333
$ws.Range("N27").Value = -1794.5625
$ws.Range("H122").Value = 4544.024
$ws.Range("I122").Value = 5797.6924
$ws.Range("J122").Value = 2506.8125
$ws.Range("K122").Value = 17393.0772
$ws.Range("L122").Value = 7520.4375
$ws.Range("M122").Value = -14943.0772
$ws.Range("N122").Value = -12420.4375
$ws.Range("H126").Value = 1928
$ws.Range("I126").Value = 1993.1538
$ws.Range("J126").Value = 1833.8889
$ws.Range("K126").Value = 5979.4614
$ws.Range("L126").Value = 5501.6667
$ws.Range("M126").Value = -3509.4614
$ws.Range("N126").Value = -10441.6667
$ws.Range("H132").Value = 2867.125
$ws.Range("I132").Value = 2157.4736
$ws.Range("J132").Value = 3904.3076
$ws.Range("K132").Value = 6472.4208
$ws.Range("L132").Value = 11712.9228
$ws.Range("M132").Value = -3942.4208
$ws.Range("N132").Value = -16772.9228
$ws.Range("H136").Value = 2498.4243
$ws.Range("I136").Value = 1189
$ws.Range("J136").Value = 5510.1
$ws.Range("K136").Value = 3567
$ws.Range("L136").Value = 16530.3
$ws.Range("M136").Value = -1017
$ws.Range("N136").Value = -21630.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6210
$ws.Range("I62").Value = 4900
$ws.Range("K62").Value = 4900
$ws.Range("M62").Value = -4276
$ws.Range("H65").Value = 6210
$ws.Range("I65").Value = 4900
$ws.Range("K65").Value = 24500
$ws.Range("M65").Value = -21380
$ws.Range("H132").Value = 2719.8
$ws.Range("I132").Value = 2438.3076
$ws.Range("J132").Value = 3242.5715
$ws.Range("K132").Value = 7314.9228
$ws.Range("L132").Value = 9727.7145
$ws.Range("M132").Value = -4784.9228
$ws.Range("N132").Value = -14787.7145
$ws.Range("H136").Value = 1545.1562
$ws.Range("I136").Value = 1482.0968
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4446.2904
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -1896.2904
$ws.Range("N136").Value = -15600
